$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" column (C2:C8) date serial value from 45185 to 45204
$ws.Range("C2:C8").Value = 45204
